$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '70.052.06'
$ws.Range("E2").Value = '  +0.14%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.588.67'
$ws.Range("E3").Value = '  -0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '579.88'
$ws.Range("E5").Value = '  -1.25%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '191.62'
$ws.Range("E6").Value = '  +0.82%  '
$ws.Range("E7").Value = '  -1.75%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.584.50'
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("E9").Value = '  +0.03%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.182'
$ws.Range("E10").Value = '  +2.25%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.664'
$ws.Range("E11").Value = '  +0.70%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '56.00'
$ws.Range("E12").Value = '  -3.38%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000305'
$ws.Range("E13").Value = '  +5.43%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.68'
$ws.Range("E14").Value = '  -0.62%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.170.12'
$ws.Range("E15").Value = '  +0.13%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '19.98'
$ws.Range("E16").Value = '  +3.53%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.590.89'
$ws.Range("E17").Value = '  -0.47%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '70.050.36'
$ws.Range("E18").Value = '  +0.23%  '
$ws.Range("E19").Value = '  +2.02%  '
$ws.Range("E20").Value = '  +0.30%  '
$ws.Range("E21").Value = '  -0.06%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '477.29'
$ws.Range("E22").Value = '  -3.37%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '19.64'
$ws.Range("E23").Value = '  +11.94%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.06'
$ws.Range("E24").Value = '  -6.19%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.40'
$ws.Range("E25").Value = '  -0.76%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '95.47'
$ws.Range("E26").Value = '  +5.44%  '
$ws.Range("E27").Value = '  -2.62%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.11'
$ws.Range("E28").Value = '  +0.30%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.46'
$ws.Range("E29").Value = '  +0.66%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '32.34'
$ws.Range("E30").Value = '  +0.44%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.70'
$ws.Range("E31").Value = '  +0.51%  '
$ws.Range("E32").Value = '  +0.25%  '
$ws.Range("E33").Value = '  +1.89%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '66.49'
$ws.Range("E34").Value = '  +2.10%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '589.27'
$ws.Range("E35").Value = '  -4.86%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '39.16'
$ws.Range("E36").Value = '  +2.78%  '
$ws.Range("E37").Value = '  +0.10%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0₃0805'
$ws.Range("E38").Value = '  -1.71%  '
$ws.Range("E39").Value = '  -1.57%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.25'
$ws.Range("E40").Value = '  +21.09%  '
$ws.Range("E41").Value = '  -5.26%  '
$ws.Range("E42").Value = '  -4.60%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.87'
$ws.Range("E43").Value = '  +7.37%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.235.12'
$ws.Range("E44").Value = '  -2.37%  '
$ws.Range("E45").Value = '  +0.64%  '
$ws.Range("E46").Value = '  -0.12%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.37'
$ws.Range("E47").Value = '  +3.13%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.46'
$ws.Range("E48").Value = '  +3.80%  '
$ws.Range("E49").Value = '  +0.79%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.999'
$ws.Range("E50").Value = '  +0.04%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.15'
$ws.Range("E51").Value = '  -5.27%  '
